# Rename the "Sequence files" worksheet to "Sequence file"
$wb = $excel.ActiveWorkbook
$seqSheet = $wb.Worksheets.Item("Sequence files")
$seqSheet.Name = "Sequence file"

# Move the active-window selection off "Collection protocol" (previously
# the active/selected tab) and onto a plain cell selection there.
$collectionSheet = $wb.Worksheets.Item("Collection protocol")
$collectionSheet.Activate()
$collectionSheet.Range("O20").Select()

# Make "Sequence file" the active/selected sheet with A4 selected,
# which also updates the workbook's activeTab.
$seqSheet.Activate()
$seqSheet.Range("A4").Select()
